# Updates the cryptocurrency price/volume table with refreshed data from the source feed.
# Rows 15/16 also swap coin name+link (WrappedBTC <-> WrappedliquidstakedEther2.0),
# reflecting a change in ranking order between the two coins.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D (Price) holds text-formatted numbers (e.g. '1.00', '0.0₃0746').
# Force the cell to text format before assigning so Excel does not silently
# reinterpret the string as a floating point number and lose exact formatting,
# then restore the default 'Normal' style so no visible formatting changes.
function Set-TextValue($cell, $text) {
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "61.310.20"
$ws.Range("E2").Value = "  +1.82%  "
Set-TextValue $ws.Range("D3") "2.358.20"
$ws.Range("E3").Value = "  +0.93%  "
Set-TextValue $ws.Range("D4") "1.00"
$ws.Range("E4").Value = "  +0.01%  "
Set-TextValue $ws.Range("D5") "558.11"
$ws.Range("E5").Value = "  +2.42%  "
Set-TextValue $ws.Range("D6") "133.07"
$ws.Range("E6").Value = "  +1.33%  "
Set-TextValue $ws.Range("D7") "1.00"
$ws.Range("E7").Value = "  +0.02%  "
Set-TextValue $ws.Range("D8") "0.586"
$ws.Range("E8").Value = "  +0.14%  "
Set-TextValue $ws.Range("D9") "2.356.41"
$ws.Range("E9").Value = "  +1.02%  "
$ws.Range("E10").Value = "  +2.27%  "
$ws.Range("E11").Value = "  +2.18%  "
$ws.Range("E12").Value = "  -0.30%  "
$ws.Range("E13").Value = "  +2.37%  "
Set-TextValue $ws.Range("D14") "24.43"
$ws.Range("E14").Value = "  +2.88%  "
$ws.Range("B15").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C15").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
Set-TextValue $ws.Range("D15") "2.778.36"
$ws.Range("E15").Value = "  +1.09%  "
$ws.Range("B16").Value = "WrappedBTC"
$ws.Range("C16").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
Set-TextValue $ws.Range("D16") "61.173.73"
$ws.Range("E16").Value = "  +1.65%  "
$ws.Range("E17").Value = "  +2.57%  "
Set-TextValue $ws.Range("D18") "2.362.91"
$ws.Range("E18").Value = "  +1.33%  "
$ws.Range("E19").Value = "  +1.84%  "
Set-TextValue $ws.Range("D20") "4.15"
$ws.Range("E20").Value = "  +0.00%  "
Set-TextValue $ws.Range("D21") "318.33"
$ws.Range("E21").Value = "  +1.48%  "
Set-TextValue $ws.Range("D22") "6.72"
$ws.Range("E22").Value = "  -0.75%  "
Set-TextValue $ws.Range("D23") "1.00"
$ws.Range("E23").Value = "  +0.29%  "
Set-TextValue $ws.Range("D24") "64.65"
$ws.Range("E24").Value = "  +1.65%  "
$ws.Range("E25").Value = "  +1.65%  "
Set-TextValue $ws.Range("D26") "0.998"
$ws.Range("E26").Value = "  -0.18%  "
Set-TextValue $ws.Range("D27") "8.11"
$ws.Range("E27").Value = "  +2.67%  "
Set-TextValue $ws.Range("D28") "1.44"
$ws.Range("E28").Value = "  +6.47%  "
$ws.Range("E29").Value = "  +10.32%  "
Set-TextValue $ws.Range("D30") "1.75"
$ws.Range("E30").Value = "  +0.64%  "
Set-TextValue $ws.Range("D31") "171.02"
$ws.Range("E31").Value = "  -0.72%  "
Set-TextValue $ws.Range("D32") "0.0₃0746"
$ws.Range("E32").Value = "  +2.35%  "
Set-TextValue $ws.Range("D33") "6.18"
$ws.Range("E33").Value = "  +4.22%  "
Set-TextValue $ws.Range("D34") "1.39"
$ws.Range("E34").Value = "  +0.52%  "
$ws.Range("E35").Value = "  +2.12%  "
Set-TextValue $ws.Range("D36") "18.23"
$ws.Range("E36").Value = "  +1.39%  "
$ws.Range("E38").Value = "  +0.03%  "
$ws.Range("E39").Value = "  +2.44%  "
Set-TextValue $ws.Range("D40") "335.40"
$ws.Range("E40").Value = "  +4.43%  "
Set-TextValue $ws.Range("D42") "38.36"
$ws.Range("E42").Value = "  +0.62%  "
Set-TextValue $ws.Range("D43") "140.50"
$ws.Range("E43").Value = "  +0.10%  "
Set-TextValue $ws.Range("D44") "3.57"
$ws.Range("E44").Value = "  +3.29%  "
Set-TextValue $ws.Range("D45") "0.0957"
$ws.Range("E45").Value = "  +1.34%  "
Set-TextValue $ws.Range("D46") "19.57"
$ws.Range("E46").Value = "  +0.65%  "
$ws.Range("E47").Value = "  +2.91%  "
Set-TextValue $ws.Range("D48") "0.0505"
$ws.Range("E48").Value = "  +1.66%  "
Set-TextValue $ws.Range("D49") "0.0₆0228"
$ws.Range("E49").Value = "  +6.84%  "
$ws.Range("E50").Value = "  +2.98%  "
Set-TextValue $ws.Range("D51") "17.47"
$ws.Range("E51").Value = "  +3.89%  "
